$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 4 ("ASSAY CATEGORY" row), shifting everything below up by one.
$ws.Rows("4:4").Delete()
